# Update Benchmark_Results.xlsx benchmark figures (2025-09-14 12:35:11 UTC)
# Applies the per-cell value changes captured in the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ŞANS OYUNLARI (ZİRAAT / G2 was blank)
$ws.Range("G2").Value = "9 TL - 9 TL"

# Row 3 - HESAPTAN EFT - Şube
$ws.Range("C3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("D3").Value = ""
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("I3").Value = ""
$ws.Range("K3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4 - HESAPTAN EFT - ATM
$ws.Range("C4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("D4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("I4").Value = ""
$ws.Range("K4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5 - HESAPTAN EFT - Mobil
$ws.Range("C5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("D5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("I5").Value = ""
$ws.Range("K5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6 - DÜZENLİ EFT
$ws.Range("C6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("D6").Value = ""
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("I6").Value = ""
$ws.Range("K6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# Row 7 - KREDİ KARTINDAN FATURA ÖDEME (ZİRAAT / G7 was blank)
$ws.Range("G7").Value = "1 TRY (Kredi kartı ile ödemelerde ek olarak nakit avans faizi uygulanır.)"

# Row 8 - HESAPTAN HAVALE - Şube
$ws.Range("C8").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("D8").Value = ""
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = "15,23 TL - 30,47 TL - 304,72 TL"
$ws.Range("I8").Value = ""
$ws.Range("K8").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 9 - HESAPTAN HAVALE - ATM
$ws.Range("C9").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("D9").Value = ""
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = "15,23 TL - 30,47 TL - 304,72 TL"
$ws.Range("I9").Value = ""
$ws.Range("K9").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 10 - HESAPTAN HAVALE - Mobil
$ws.Range("C10").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("D10").Value = ""
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = "15,23 TL - 30,47 TL - 304,72 TL"
$ws.Range("I10").Value = ""
$ws.Range("K10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 11 - DÜZENLİ HAVALE
$ws.Range("C11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("D11").Value = ""
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = "3,05 TL - 6,1 TL - 76,18 TL"
$ws.Range("I11").Value = ""
$ws.Range("K11").Value = "3,05 TL - 6,09 TL - 76,17 TL"

# Row 12 - GİDEN SWIFT
$ws.Range("C12").Value = "WU: 1.000,01 USD–9,51 USD"
$ws.Range("D12").Value = ""
$ws.Range("G12").Value = ""
$ws.Range("K12").Value = "WU: ,USD–; Diğer: 404,16 TL–3.403,42 TL"

# Row 13 - GELEN SWIFT
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("D13").Value = ""
$ws.Range("H13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"
$ws.Range("I13").Value = ""
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# Row 14 - GİDEN SWIFT - Mobil
$ws.Range("C14").Value = "40.000 TL - 1.904,76 TL"
$ws.Range("D14").Value = ""
$ws.Range("G14").Value = ""
$ws.Range("H14").Value = "2.100 TL - 4.300 TL"
$ws.Range("K14").Value = "914,14 TL - 4.265,98 TL"

# Row 24 - SENET TAHSİLE ALMA
$ws.Range("C24").Value = "457,14 TL"

# Row 25 - MUAMELESİZ SENET İADESİ
$ws.Range("C25").Value = "380,95 TL"
